$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CZZ")

# Updated financial figures (yearly financials refresh)
# Row 8
$ws.Range("D8").Value = 3482300
$ws.Range("E8").Value = 3209400
$ws.Range("F8").Value = 3167700
$ws.Range("G8").Value = 2283000
$ws.Range("H8").Value = 1763400
$ws.Range("I8").Value = 1175800
$ws.Range("J8").Value = 1170000

# Row 9
$ws.Range("D9").Value = 2367000
$ws.Range("E9").Value = 2132400
$ws.Range("F9").Value = 2216600
$ws.Range("G9").Value = 1628800
$ws.Range("H9").Value = 1250700
$ws.Range("I9").Value = 823300
$ws.Range("J9").Value = 947600

# Row 10
$ws.Range("D10").Value = 1115300
$ws.Range("E10").Value = 1077000
$ws.Range("F10").Value = 951100
$ws.Range("G10").Value = 654200
$ws.Range("H10").Value = 512800
$ws.Range("I10").Value = 352500
$ws.Range("J10").Value = 222400

# Row 14
$ws.Range("D14").Value = 44300
$ws.Range("E14").Value = 28700
$ws.Range("F14").Value = 30100
$ws.Range("G14").Value = 40400
$ws.Range("H14").Value = 22100
$ws.Range("I14").Value = -18800
$ws.Range("J14").Value = -679900

# Row 17
$ws.Range("D17").Value = 2398800
$ws.Range("E17").Value = 2283400
$ws.Range("F17").Value = 2436300
$ws.Range("G17").Value = 1907300
$ws.Range("H17").Value = 1505700
$ws.Range("I17").Value = 1004100
$ws.Range("J17").Value = 399500

# Row 18
$ws.Range("D18").Value = 1083500
$ws.Range("E18").Value = 926000
$ws.Range("F18").Value = 731400
$ws.Range("G18").Value = 375700
$ws.Range("H18").Value = 257800
$ws.Range("I18").Value = 171700
$ws.Range("J18").Value = 770500

# Row 20
$ws.Range("D20").Value = -247200
$ws.Range("E20").Value = -312900
$ws.Range("F20").Value = 130700
$ws.Range("G20").Value = -96700
$ws.Range("H20").Value = 66300
$ws.Range("I20").Value = 184800
$ws.Range("J20").Value = 175200

# Row 21
$ws.Range("D21").Value = 1321700
$ws.Range("E21").Value = 1047700
$ws.Range("F21").Value = 1157200
$ws.Range("G21").Value = 448800
$ws.Range("H21").Value = 434000
$ws.Range("I21").Value = 440400
$ws.Range("J21").Value = 1011100

# Row 22
$ws.Range("D22").Value = 458200
$ws.Range("E22").Value = 470600
$ws.Range("F22").Value = 690800
$ws.Range("G22").Value = 158800
$ws.Range("H22").Value = 185800
$ws.Range("I22").Value = 120400
$ws.Range("J22").Value = 142000

# Row 23
$ws.Range("D23").Value = 378100
$ws.Range("E23").Value = 142500
$ws.Range("F23").Value = 171400
$ws.Range("G23").Value = 120100
$ws.Range("H23").Value = 138200
$ws.Range("I23").Value = 236100
$ws.Range("J23").Value = 803700

# Row 24
$ws.Range("D24").Value = 109800
$ws.Range("E24").Value = 15800
$ws.Range("F24").Value = -7800
$ws.Range("G24").Value = 6000
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 32400
$ws.Range("J24").Value = 261300

# Row 26
$ws.Range("D26").Value = 268200
$ws.Range("E26").Value = 126700
$ws.Range("F26").Value = 179200
$ws.Range("G26").Value = 114100
$ws.Range("H26").Value = 128200
$ws.Range("I26").Value = 203700
$ws.Range("J26").Value = 542500

# Row 27
$ws.Range("D27").Value = 141300
$ws.Range("E27").Value = 80300
$ws.Range("F27").Value = 78200
$ws.Range("G27").Value = -5200
$ws.Range("H27").Value = 31400
$ws.Range("I27").Value = 73100
$ws.Range("J27").Value = 285100

# Row 29
$ws.Range("E29").Value = 9000
$ws.Range("F29").Value = 25900
$ws.Range("G29").Value = 46300
$ws.Range("I29").Value = 35600
$ws.Range("J29").Value = 16500

# Row 32
$ws.Range("D32").Value = 247200
$ws.Range("E32").Value = 312900
$ws.Range("F32").Value = -130700
$ws.Range("G32").Value = 96700
$ws.Range("H32").Value = -66300
$ws.Range("I32").Value = -184800
$ws.Range("J32").Value = -175200

# Row 33
$ws.Range("D33").Value = 141300
$ws.Range("E33").Value = 89300
$ws.Range("F33").Value = 104000
$ws.Range("G33").Value = 41100
$ws.Range("H33").Value = 31400
$ws.Range("I33").Value = 108700
$ws.Range("J33").Value = 301600

# Row 35
$ws.Range("D35").Value = 141300
$ws.Range("E35").Value = 89300
$ws.Range("F35").Value = 104000
$ws.Range("G35").Value = 41100
$ws.Range("H35").Value = 31400
$ws.Range("I35").Value = 108700
$ws.Range("J35").Value = 301600

# Row 41
$ws.Range("D41").Value = 13700
$ws.Range("E41").Value = 96400
$ws.Range("F41").Value = 78300
$ws.Range("G41").Value = 462000
$ws.Range("H41").Value = 31500
$ws.Range("I41").Value = 682200
$ws.Range("J41").Value = 424100

# Row 42
$ws.Range("D42").Value = 2485600
$ws.Range("E42").Value = 1388300
$ws.Range("F42").Value = 975700
$ws.Range("G42").Value = 422100
$ws.Range("H42").Value = 378100
$ws.Range("I42").Value = 407200

# Row 43
$ws.Range("D43").Value = 502800
$ws.Range("E43").Value = 481300
$ws.Range("F43").Value = 368800
$ws.Range("G43").Value = 303100
$ws.Range("H43").Value = 287800
$ws.Range("I43").Value = 856900
$ws.Range("J43").Value = 499800

# Row 44
$ws.Range("D44").Value = 170000
$ws.Range("E44").Value = 161700
$ws.Range("F44").Value = 336800
$ws.Range("G44").Value = 181400
$ws.Range("H44").Value = 80000
$ws.Range("I44").Value = 257500
$ws.Range("J44").Value = 191800

# Row 45
$ws.Range("D45").Value = 169600
$ws.Range("E45").Value = 120500
$ws.Range("F45").Value = 161100
$ws.Range("G45").Value = 59000
$ws.Range("H45").Value = 136400
$ws.Range("I45").Value = 257900
$ws.Range("J45").Value = 103000

# Row 46
$ws.Range("D46").Value = 3341600
$ws.Range("E46").Value = 2248200
$ws.Range("F46").Value = 1752400
$ws.Range("G46").Value = 914000
$ws.Range("H46").Value = 913700
$ws.Range("I46").Value = 882800
$ws.Range("J46").Value = 1218700

# Row 47
$ws.Range("D47").Value = 2561300
$ws.Range("E47").Value = 2536400
$ws.Range("F47").Value = 2473000
$ws.Range("G47").Value = 2472700
$ws.Range("H47").Value = 2517500
$ws.Range("I47").Value = 3465700
$ws.Range("J47").Value = 480700

# Row 48
$ws.Range("D48").Value = 2994900
$ws.Range("E48").Value = 2750000
$ws.Range("F48").Value = 3179300
$ws.Range("G48").Value = 1413600
$ws.Range("H48").Value = 911000
$ws.Range("I48").Value = 1030700
$ws.Range("J48").Value = 2342800

# Row 49
$ws.Range("D49").Value = 4351700
$ws.Range("E49").Value = 4386500
$ws.Range("F49").Value = 4437900
$ws.Range("G49").Value = 5274400
$ws.Range("H49").Value = 2583800
$ws.Range("I49").Value = 5839500
$ws.Range("J49").Value = 1264500

# Row 52
$ws.Range("D52").Value = 1011400
$ws.Range("E52").Value = 1018300
$ws.Range("F52").Value = 1553100
$ws.Range("G52").Value = 549100
$ws.Range("H52").Value = 410500
$ws.Range("I52").Value = 594600
$ws.Range("J52").Value = 454400

# Row 54
$ws.Range("D54").Value = 14261000
$ws.Range("E54").Value = 12939500
$ws.Range("F54").Value = 13395700
$ws.Range("G54").Value = 7618500
$ws.Range("H54").Value = 7336500
$ws.Range("I54").Value = 7055800
$ws.Range("J54").Value = 5683500

# Row 57
$ws.Range("D57").Value = 624000
$ws.Range("E57").Value = 521100
$ws.Range("F57").Value = 503500
$ws.Range("G57").Value = 285200
$ws.Range("H57").Value = 221100
$ws.Range("I57").Value = 235100
$ws.Range("J57").Value = 200400

# Row 58
$ws.Range("D58").Value = 1074800
$ws.Range("E58").Value = 744600
$ws.Range("F58").Value = 993500
$ws.Range("G58").Value = 270800
$ws.Range("H58").Value = 269400
$ws.Range("I58").Value = 412400
$ws.Range("J58").Value = 138500

# Row 59
$ws.Range("D59").Value = 614300
$ws.Range("E59").Value = 433900
$ws.Range("F59").Value = 499500
$ws.Range("G59").Value = 205700
$ws.Range("H59").Value = 189100
$ws.Range("I59").Value = 379600
$ws.Range("J59").Value = 193000

# Row 60
$ws.Range("D60").Value = 2313100
$ws.Range("E60").Value = 1699600
$ws.Range("F60").Value = 1774800
$ws.Range("G60").Value = 761700
$ws.Range("H60").Value = 679600
$ws.Range("I60").Value = 765800
$ws.Range("J60").Value = 531900

# Row 61
$ws.Range("D61").Value = 5849800
$ws.Range("E61").Value = 5437600
$ws.Range("F61").Value = 5512900
$ws.Range("G61").Value = 2403100
$ws.Range("H61").Value = 2061800
$ws.Range("I61").Value = 1768900
$ws.Range("J61").Value = 1194500

# Row 62
$ws.Range("D62").Value = 1724400
$ws.Range("E62").Value = 1697800
$ws.Range("F62").Value = 2142700
$ws.Range("G62").Value = 1015600
$ws.Range("H62").Value = 1180200
$ws.Range("I62").Value = 2048800
$ws.Range("J62").Value = 1526200

# Row 66
$ws.Range("D66").Value = 12712800
$ws.Range("E66").Value = 11331300
$ws.Range("F66").Value = 11879500
$ws.Range("G66").Value = 6132700
$ws.Range("H66").Value = 5827400
$ws.Range("I66").Value = 5518200
$ws.Range("J66").Value = 4253500

# Row 72
$ws.Range("D72").Value = 815800
$ws.Range("E72").Value = 691200
$ws.Range("F72").Value = 610200
$ws.Range("G72").Value = 530300
$ws.Range("H72").Value = 547900
$ws.Range("I72").Value = 1570600
$ws.Range("J72").Value = 1473800

# Row 76
$ws.Range("D76").Value = 1548200
$ws.Range("E76").Value = 1608100
$ws.Range("F76").Value = 1516200
$ws.Range("G76").Value = 1485800
$ws.Range("H76").Value = 1509100
$ws.Range("I76").Value = 1537600
$ws.Range("J76").Value = 1429900

# Row 81
$ws.Range("D81").Value = 141300
$ws.Range("E81").Value = 89300
$ws.Range("F81").Value = 104000
$ws.Range("G81").Value = 41100
$ws.Range("H81").Value = 31400
$ws.Range("I81").Value = 108700
$ws.Range("J81").Value = 301600

# Row 83
$ws.Range("D83").Value = 497000
$ws.Range("E83").Value = 444900
$ws.Range("F83").Value = 302000
$ws.Range("G83").Value = 173900
$ws.Range("H83").Value = 112600
$ws.Range("I83").Value = 85800
$ws.Range("J83").Value = 66900

# Row 89
$ws.Range("D89").Value = 1048100
$ws.Range("E89").Value = 932000
$ws.Range("F89").Value = 859000
$ws.Range("G89").Value = 286400
$ws.Range("H89").Value = 289800
$ws.Range("I89").Value = 72600
$ws.Range("J89").Value = 167200

# Row 91
$ws.Range("D91").Value = -630000
$ws.Range("E91").Value = -561700
$ws.Range("F91").Value = -515400
$ws.Range("G91").Value = -272300
$ws.Range("H91").Value = -250300
$ws.Range("I91").Value = -975500
$ws.Range("J91").Value = -233200

# Row 94
$ws.Range("D94").Value = -917200
$ws.Range("E94").Value = -186400
$ws.Range("F94").Value = -257200
$ws.Range("G94").Value = -79100
$ws.Range("H94").Value = -129400
$ws.Range("I94").Value = -694900
$ws.Range("J94").Value = -302600

# Row 96
$ws.Range("D96").Value = -279400
$ws.Range("E96").Value = -369300
$ws.Range("F96").Value = -168400
$ws.Range("G96").Value = -152800
$ws.Range("H96").Value = -80600
$ws.Range("I96").Value = -101600
$ws.Range("J96").Value = -85500

# Row 100
$ws.Range("D100").Value = -145000
$ws.Range("E100").Value = -466400
$ws.Range("F100").Value = -139100
$ws.Range("G100").Value = -168800
$ws.Range("H100").Value = -172600
$ws.Range("I100").Value = 751500
$ws.Range("J100").Value = 71000

# Row 101
$ws.Range("D101").Value = 28400
$ws.Range("E101").Value = -24400
$ws.Range("F101").Value = 13200
$ws.Range("G101").Value = -2600
$ws.Range("H101").Value = 3500
$ws.Range("I101").Value = 1000
$ws.Range("J101").Value = 4000

# Row 102
$ws.Range("D102").Value = 14300
$ws.Range("E102").Value = 254800
$ws.Range("F102").Value = 476000
$ws.Range("G102").Value = 35800
$ws.Range("H102").Value = -8800
$ws.Range("I102").Value = 130200
$ws.Range("J102").Value = -60400
